$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$NewValue)
    $Cell.NumberFormat = "@"
    $Cell.Value = $NewValue
    $Cell.NumberFormat = "General"
}

Set-TextValue $ws.Range("D2") "69.127.13"
Set-TextValue $ws.Range("E2") "  -2.77%  "
Set-TextValue $ws.Range("D3") "3.676.54"
Set-TextValue $ws.Range("E3") "  -4.10%  "
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  +0.07%  "
Set-TextValue $ws.Range("D5") "677.80"
Set-TextValue $ws.Range("E5") "  -4.07%  "
Set-TextValue $ws.Range("E6") "  -5.96%  "
Set-TextValue $ws.Range("D7") "3.676.90"
Set-TextValue $ws.Range("E7") "  -4.04%  "
Set-TextValue $ws.Range("E8") "  -0.06%  "
Set-TextValue $ws.Range("D9") "0.487"
Set-TextValue $ws.Range("E9") "  -7.04%  "
Set-TextValue $ws.Range("D10") "0.147"
Set-TextValue $ws.Range("E10") "  -8.79%  "
Set-TextValue $ws.Range("D11") "7.24"
Set-TextValue $ws.Range("E11") "  -2.53%  "
Set-TextValue $ws.Range("D12") "0.446"
Set-TextValue $ws.Range("E12") "  -2.59%  "
Set-TextValue $ws.Range("D13") "0.0000231"
Set-TextValue $ws.Range("E13") "  -9.26%  "
Set-TextValue $ws.Range("E14") "  -9.54%  "
Set-TextValue $ws.Range("D15") "4.296.41"
Set-TextValue $ws.Range("E15") "  -4.12%  "
Set-TextValue $ws.Range("D16") "3.672.44"
Set-TextValue $ws.Range("E16") "  -2.67%  "
Set-TextValue $ws.Range("D17") "69.156.28"
Set-TextValue $ws.Range("E17") "  -2.80%  "
Set-TextValue $ws.Range("E18") "  -1.68%  "
Set-TextValue $ws.Range("D19") "16.16"
Set-TextValue $ws.Range("E19") "  -6.98%  "
Set-TextValue $ws.Range("D20") "6.52"
Set-TextValue $ws.Range("E20") "  -9.55%  "
Set-TextValue $ws.Range("D21") "480.48"
Set-TextValue $ws.Range("E21") "  -2.99%  "
Set-TextValue $ws.Range("D22") "9.76"
Set-TextValue $ws.Range("E22") "  -8.18%  "
Set-TextValue $ws.Range("E23") "  -9.85%  "
Set-TextValue $ws.Range("D24") "78.87"
Set-TextValue $ws.Range("E24") "  -7.91%  "
Set-TextValue $ws.Range("D25") "3.818.53"
Set-TextValue $ws.Range("E25") "  -4.24%  "
Set-TextValue $ws.Range("D26") "11.55"
Set-TextValue $ws.Range("E26") "  -4.63%  "
Set-TextValue $ws.Range("E27") "  -0.17%  "
Set-TextValue $ws.Range("E28") "  -12.67%  "
Set-TextValue $ws.Range("E29") "  -12.01%  "
Set-TextValue $ws.Range("D30") "1.81"
Set-TextValue $ws.Range("E30") "  -13.10%  "
Set-TextValue $ws.Range("D31") "2.71"
Set-TextValue $ws.Range("E31") "  -12.41%  "
Set-TextValue $ws.Range("D32") "2.09"
Set-TextValue $ws.Range("E32") "  -6.26%  "
Set-TextValue $ws.Range("E33") "  -10.27%  "
Set-TextValue $ws.Range("D34") "0.167"
Set-TextValue $ws.Range("E34") "  -4.67%  "
Set-TextValue $ws.Range("D35") "0.998"
Set-TextValue $ws.Range("E35") "  -0.17%  "
Set-TextValue $ws.Range("D36") "26.61"
Set-TextValue $ws.Range("E36") "  -9.27%  "
Set-TextValue $ws.Range("D37") "3.641.53"
Set-TextValue $ws.Range("E37") "  -4.24%  "
Set-TextValue $ws.Range("E38") "  -7.53%  "
Set-TextValue $ws.Range("D39") "6.03"
Set-TextValue $ws.Range("E39") "  +0.96%  "
Set-TextValue $ws.Range("D40") "0.0926"
Set-TextValue $ws.Range("E40") "  -9.35%  "
Set-TextValue $ws.Range("D42") "2.17"
Set-TextValue $ws.Range("E42") "  -6.36%  "
Set-TextValue $ws.Range("E43") "  +0.07%  "
Set-TextValue $ws.Range("D44") "0.947"
Set-TextValue $ws.Range("E44") "  -9.75%  "
Set-TextValue $ws.Range("D45") "159.68"
Set-TextValue $ws.Range("E45") "  -2.62%  "
Set-TextValue $ws.Range("E46") "  -2.27%  "
Set-TextValue $ws.Range("D47") "2.85"
Set-TextValue $ws.Range("E47") "  -14.37%  "
Set-TextValue $ws.Range("E48") "  -4.00%  "
$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D49") "384.97"
Set-TextValue $ws.Range("E49") "  -10.45%  "
$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
Set-TextValue $ws.Range("D50") "0.000272"
Set-TextValue $ws.Range("E50") "  -12.75%  "
Set-TextValue $ws.Range("E51") "  -8.96%  "
